$d = $word.ActiveDocument

# Locate the final paragraph of the document: "The problem of "small group
# classification disadvantages" is alleviated."
$lastIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastIndex)
$targetRange = $targetPara.Range
$targetRange.Collapse(0)

# 1) Insert a blank paragraph right after it.
$targetRange.InsertParagraphAfter()

# 2) Insert the "Thanks to my teammate..." paragraph.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("Thanks to my teammate Shakin Shahria for giving me the floor. After reading the paper thoroughly I found some findings in this paper. First of all, the authors evaluated the importance of different features  using the information gain rate. As IPv6 is much larger, it has more features and using many features will affect the performance. So to enhance the classification efficiency the algorithm needs to choose the essential features.")

# 3) Insert the "Secondly from the testing..." paragraph.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("Secondly from the testing we saw that the new Information Gain Ratio Average Distance KNN algorithm performs better than the already existing GR-KNN and TAD-KNN.")

# 4) Insert the "Thirdly, this new algorithm..." paragraph.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("Thirdly, this new algorithm improves the IPv6 DDoS attack classification accuracy, algorithm stability and solve the problems on distant sample effect and small group classification of traditional KNN algorithm.")
